# TestSpreadSheet.xlsx update:
#  - FindCarTest / CarNameAndPrice sheets: result column (C3) changes from "y" to "n"
#  - CarNameAndPrice's selection (B5) is cleared and both sheets end up with C4 selected
#  - FindCarTest becomes the active (visible) tab instead of CarNameAndPrice

$wb = $excel.ActiveWorkbook

$findCarTest      = $wb.Worksheets.Item("FindCarTest")
$carNameAndPrice  = $wb.Worksheets.Item("CarNameAndPrice")

# Update the "found it?" result cell on both sheets from "y" to "n"
$findCarTest.Range("C3").Value = "n"
$carNameAndPrice.Range("C3").Value = "n"

# CarNameAndPrice is no longer the active sheet; move its selection to C4
$carNameAndPrice.Range("C4").Select()

# FindCarTest becomes the active sheet, selected at C4
$findCarTest.Activate()
$findCarTest.Range("C4").Select()
